$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 539, shifting rows 539:577 down to 540:578
$ws.Rows("539:539").Insert()

# Populate the newly inserted row 539 with the weekly data point
$ws.Range("A539").Value = 9
$ws.Range("B539").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C539").Value = "Metropolitana"
$ws.Range("D539").Value = 45106
$ws.Range("E539").Value = 13
$ws.Range("F539").Value = 100112044
$ws.Range("G539").Value = "Perejil"
$ws.Range("H539").Value = "Sin especificar"
$ws.Range("I539").Value = "Primera"
$ws.Range("J539").Value = 52
$ws.Range("K539").Value = 16000
$ws.Range("L539").Value = 18000
$ws.Range("M539").Value = 17000
$ws.Range("N539").Value = "$/docena de atados"
$ws.Range("O539").Value = "Región Metropolitana"
$ws.Range("P539").Value = 5667
$ws.Range("Q539").Value = 3
$ws.Range("R539").Value = "Hortaliza"
